$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- Header row (row 1): extend with new columns G..M, copying the
# existing header-cell format (style) so the new cells match the
# look of the existing headers (bold font, border, centered). ---
$ws.Cells.Item(1, 2).Copy() | Out-Null
$ws.Range($ws.Cells.Item(1, 7), $ws.Cells.Item(1, 13)).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(1, 2).Value = "bank"
$ws.Cells.Item(1, 3).Value = "deposit_type"
$ws.Cells.Item(1, 4).Value = "currency"
$ws.Cells.Item(1, 5).Value = "owner"
$ws.Cells.Item(1, 6).Value = "total"
$ws.Cells.Item(1, 7).Value = "property_category"
$ws.Cells.Item(1, 8).Value = "category"
$ws.Cells.Item(1, 9).Value = "date"
$ws.Cells.Item(1, 10).Value = "legislator_name"
$ws.Cells.Item(1, 11).Value = "legislator_id"
$ws.Cells.Item(1, 12).Value = "source_file"
$ws.Cells.Item(1, 13).Value = "index"

# --- Data row (row 2): extend with new columns G..M, copying the
# existing data-cell format (style) for consistency. ---
$ws.Cells.Item(2, 2).Copy() | Out-Null
$ws.Range($ws.Cells.Item(2, 7), $ws.Cells.Item(2, 13)).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(2, 2).Value = "臺灣銀行群賢分行"
$ws.Cells.Item(2, 3).Value = "活期儲蓄存款"
$ws.Cells.Item(2, 4).Value = "新臺幣"
$ws.Cells.Item(2, 5).Value = "林正二"
$ws.Cells.Item(2, 6).Value = 673238
$ws.Cells.Item(2, 7).Value = "deposit"
$ws.Cells.Item(2, 8).Value = "normal"

# "date" (I2) must stay plain text ("2012-04-12"), not an auto-converted
# date serial. Force text format before assigning, then reapply the
# regular data-cell format (which doesn't touch the already-stored value).
$ws.Cells.Item(2, 9).NumberFormat = "@"
$ws.Cells.Item(2, 9).Value = "2012-04-12"
$ws.Cells.Item(2, 2).Copy() | Out-Null
$ws.Cells.Item(2, 9).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(2, 10).Value = "林正二"
$ws.Cells.Item(2, 11).Value = 788
$ws.Cells.Item(2, 12).Value = "tmp32921"
$ws.Cells.Item(2, 13).Value = 53
